$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Recomputed NATMI LR-pair statistics using updated TPM-based expression values
# (commit: "update scripts wuth new tpm"). Columns E (col 5) through T (col 20).
$data = @{
  2 = @(3,1,5.898379333333334,17.695138,0.7828421971437433,0.7828421971437433,2,0.6666666666666666,0.07259900000000001,0.217797,0.0162094769588191,0.0162094769588191,0.4282164412206668,3.853947970986,0.01268946255699283,0.01268946255699283)
  3 = @(3,1,5.898379333333334,17.695138,0.7828421971437433,0.7828421971437433,3,1,4.187332,12.561996,0.9349228167457665,0.9349228167457664,24.69847253060534,222.286252775448,0.731897032021073,0.731897032021073)
  4 = @(3,1,5.898379333333334,17.695138,0.7828421971437433,0.7828421971437433,3,1,0.2188686666666667,0.656606,0.04886770629541442,0.04886770629541441,1.290970420180889,11.618733781628,0.03825570256567735,0.03825570256567735)
  5 = @(2,0.6666666666666666,0.7963256666666666,2.388977,0.1056895969732402,0.1056895969732402,2,0.6666666666666666,0.07259900000000001,0.217797,0.0162094769588191,0.0162094769588191,0.05781244707433333,0.520312023669,0.001713173086924615,0.001713173086924614)
  6 = @(2,0.6666666666666666,0.7963256666666666,2.388977,0.1056895969732402,0.1056895969732402,3,1,4.187332,12.561996,0.9349228167457665,0.9349228167457664,3.334479946454667,30.010319518092,0.0988116157029466,0.09881161570294657)
  7 = @(2,0.6666666666666666,0.7963256666666666,2.388977,0.1056895969732402,0.1056895969732402,3,1,0.2188686666666667,0.656606,0.04886770629541442,0.04886770629541441,0.1742907368957778,1.568616632062,0.005164808183369025,0.005164808183369023)
  8 = @(1,0.3333333333333333,0.1290883333333333,0.387265,0.01713280695956549,0.01713280695956549,2,0.6666666666666666,0.07259900000000001,0.217797,0.0162094769588191,0.0162094769588191,0.009371683911666668,0.08434515520500001,0.0002809641680332896,0.0002777138396509724)
  9 = @(1,0.3333333333333333,0.1290883333333333,0.387265,0.01713280695956549,0.01713280695956549,3,1,4.187332,12.561996,0.9349228167457665,0.9349228167457664,0.5405357089933334,4.86482138094,0.01601785214139844,0.01601785214139844)
  10 = @(1,0.3333333333333333,0.1290883333333333,0.387265,0.01713280695956549,0.01713280695956549,3,1,0.2188686666666667,0.656606,0.04886770629541442,0.04886770629541441,0.02825339139888889,0.25428052259,0.0008372409785160783,0.0008372409785160783)
  11 = @(3,1,0.4975636666666667,1.492691,0.06603743367792279,0.06603743367792278,2,0.6666666666666666,0.07259900000000001,0.217797,0.0162094769588191,0.0162094769588191,0.03612262463633334,0.325103621727,0.001070432259621834,0.001070432259621834)
  12 = @(3,1,0.4975636666666667,1.492691,0.06603743367792279,0.06603743367792278,3,1,4.187332,12.561996,0.9349228167457665,0.9349228167457664,2.083464263470667,18.751178371236,0.06173990350482532,0.0617399035048253)
  13 = @(3,1,0.4975636666666667,1.492691,0.06603743367792279,0.06603743367792278,3,1,0.2188686666666667,0.656606,0.04886770629541442,0.04886770629541441,0.1089010963051111,0.980109866746,0.00322709791347564,0.003227097913475639)
  14 = @(2,0.6666666666666666,0.213213,0.639639,0.02829796524552828,0.02829796524552828,2,0.6666666666666666,0.07259900000000001,0.217797,0.0162094769588191,0.0162094769588191,0.015479050587,0.139311455283,0.0004586952156288544,0.0004586952156288543)
  15 = @(2,0.6666666666666666,0.213213,0.639639,0.02829796524552828,0.02829796524552828,3,1,4.187332,12.561996,0.9349228167457665,0.9349228167457664,0.8927936177160001,8.035142559443999,0.0264564133755231,0.0264564133755231)
  16 = @(2,0.6666666666666666,0.213213,0.639639,0.02829796524552828,0.02829796524552828,3,1,0.2188686666666667,0.656606,0.04886770629541442,0.04886770629541441,0.046665645026,0.419990805234,0.001382856654376321,0.00138285665437632)
}

foreach ($r in $data.Keys) {
    $rowVals = $data[$r]
    for ($j = 0; $j -lt $rowVals.Count; $j++) {
        $ws.Cells.Item($r, 5 + $j).Value = $rowVals[$j]
    }
}